# Actualizacion SmartScore desde Streamlit (Juan)
# - Convierte a numero los SmartScore de la fila 2 (Chava), que se habian
#   guardado como texto.
# - Agrega la fila 3 con la nueva sesion de Juan (misma estructura de columnas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 2 (Chava): normalizar las columnas de SmartScore a tipo numerico ---
$ws.Range("E2").Value = 0.548
$ws.Range("H2").Value = 0.489
$ws.Range("K2").Value = 0.477
$ws.Range("N2").Value = 0.605
$ws.Range("Q2").Value = 0.602
$ws.Range("T2").Value = 0.555
$ws.Range("W2").Value = 0.719
$ws.Range("Z2").Value = 0.601
$ws.Range("AC2").Value = 0.576

# --- Fila 3 (Juan): nueva fila de resultados ---
$ws.Range("A3").Value = "Juan"
$ws.Range("B3").Value = "2025-10-28 02:50:18"

$pesosJuan = @"
{
  "portion": 0.4,
  "diet": 0.8571428571428571,
  "salt": 0.6,
  "fat": 1.0,
  "natural": 0.8,
  "convenience": 0.2,
  "price": 0.8
}
"@
$ws.Range("C3").Value = $pesosJuan

# Los puntajes SmartScore de esta fila llegaron como texto desde Streamlit,
# igual que en la exportacion original, asi que se conservan como texto.
$ws.Range("D3").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("E3").Value = "'0.533"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("G3").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("H3").Value = "'0.475"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("J3").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("K3").Value = "'0.426"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("M3").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("N3").Value = "'0.703"
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("P3").Value = "Annie’s Shells & White Cheddar"
$ws.Range("Q3").Value = "'0.639"
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("S3").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("T3").Value = "'0.552"
$ws.Range("T3").Style = "Normal"
$ws.Range("U3").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("V3").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("W3").Value = "'0.698"
$ws.Range("W3").Style = "Normal"
$ws.Range("X3").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("Y3").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("Z3").Value = "'0.601"
$ws.Range("Z3").Style = "Normal"
$ws.Range("AA3").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AB3").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AC3").Value = "'0.579"
$ws.Range("AC3").Style = "Normal"
$ws.Range("AD3").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# La celda C3 tiene saltos de linea; volver a autoajustar el alto de fila
# para que no quede con un alto "personalizado" distinto al resto.
$ws.Rows.Item(3).AutoFit()

Write-Output "Fila 3 (Juan) agregada; SmartScores de fila 2 convertidos a numero."
